$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last fully-styled data row (row 23) down onto
# the two new data rows so the new cells pick up the same number formats /
# alignment as the rest of the table (A/B integer columns, C/D amount
# columns, E date column).
$ws.Range("A23:E23").Copy() | Out-Null
$ws.Range("A25:E26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New topup rows for 1-Feb-2021 and 2-Feb-2021 (order ids 26155067 / 26165120).
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 26155067
$ws.Range("C25").Value = 73913
$ws.Range("D25").Value = 70000.45
$ws.Range("E25").Value = 44228

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 26165120
$ws.Range("C26").Value = 38012
$ws.Range("D26").Value = 35999.4
$ws.Range("E26").Value = 44229

# The F:I columns already carry shared formulas (=IF(B="","",...)) all the
# way down to row 220; re-assert them on the two rows we just populated so
# they pick up the now-non-blank B25/B26 and recompute instead of keeping
# their previously-cached "" result.
$ws.Range("F25:I26").Formula = $ws.Range("F25:I26").Formula

# Re-point the frozen-pane viewport / active selection the way it was left
# after entering the new rows.
$ws.Range("A9").Select()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("H27").Select()
